$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (CT-23): Status Passed -> Failed
$ws.Range("E7").Value = "Failed"
$ws.Range("E7").Style = "Bad"

# Row 22 (was CT-96, before insertion shift): Test Case ID CT-96 -> CT-100
$ws.Range("B22").Value = "CT-100"

# Insert a new row before row 14 (old CT-46), shifting everything below it down
$ws.Rows("14").Insert()

# Fill in the new row 14: CT-44, with raw serial date values, Failed status
$ws.Range("B14").Value = "CT-44"
$ws.Range("C14").Value = 44540.448611111111
$ws.Range("D14").Value = 44540.490277777797
$ws.Range("D14").NumberFormat = "[h]:mm:ss;@"
$ws.Range("E14").Value = "Failed"

# Row 12 (was CT-28): Test Case ID -> CL-29
$ws.Range("B12").Value = "CL-29"

# Row 16 (CT-58, after insertion shift): Status Passed -> Failed
$ws.Range("E16").Value = "Failed"
$ws.Range("E16").Style = "Bad"

# Row 20 (was CT-71, after shift): Test Case ID CT-71 -> CT-73
$ws.Range("B20").Value = "CT-73"

# Row 22 (was CT-74, after shift): Status Passed -> Failed
$ws.Range("E22").Value = "Failed"
$ws.Range("E22").Style = "Bad"

# Final selection as left by the editor
$ws.Range("B12").Select() | Out-Null
